# Rename header cells: CLAVE_* -> ID_* (see commit: "arreglar cambiar clave por id")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demanda")

$ws.Range("B1").Value = "ID_BARRIO"
$ws.Range("D1").Value = "ID_MUNICIPIO"
$ws.Range("E1").Value = "ID_SECTOR"
$ws.Range("F1").Value = "ID_CORREGIMIENTO"

# Update the active selection to match the saved view state (B2)
$ws.Range("B2").Select()
